$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.097.20"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.638.94"
$ws.Range("E3").Value = "  -1.88%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.37"
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5262"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06323"
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.78"
$ws.Range("E10").Value = "  -1.93%  "
$ws.Range("E11").Value = "  +1.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.631.24"
$ws.Range("E12").Value = "  -2.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.431"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.861.04"
$ws.Range("E14").Value = "  -2.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5514"
$ws.Range("E15").Value = "  -1.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8178"
$ws.Range("E16").Value = "  +2.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.24"
$ws.Range("E17").Value = "  -2.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.072.98"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.700"
$ws.Range("E20").Value = "  -1.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "188.81"
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.16"
$ws.Range("E22").Value = "  -2.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.171"
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.03"
$ws.Range("E25").Value = "  -2.37%  "
$ws.Range("E26").Value = "  -2.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.410"
$ws.Range("E27").Value = "  -1.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.89"
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.407"
$ws.Range("E29").Value = "  +3.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05996"
$ws.Range("E30").Value = "  -4.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.257"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.452"
$ws.Range("E32").Value = "  -1.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.416"
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.644"
$ws.Range("E34").Value = "  +0.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9889"
$ws.Range("E35").Value = "  -1.14%  "
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.766"
$ws.Range("E36").Value = "  +1.07%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.397"
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5742"
$ws.Range("E38").Value = "  -5.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01620"
$ws.Range("E39").Value = "  +0.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8557"
$ws.Range("E40").Value = "  -2.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.041.87"
$ws.Range("E41").Value = "  -5.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.720"
$ws.Range("E43").Value = "  -6.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.61"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.788.26"
$ws.Range("E45").Value = "  -1.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈109"
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.57"
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9991"
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.061"
$ws.Range("E49").Value = "  +0.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05171"
$ws.Range("E50").Value = "  -1.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4220"
$ws.Range("E51").Value = "  -0.62%  "
